$d = $word.ActiveDocument

# 1. "Варіант 14" -> "Варіант 10"
$rng = $d.Content
$found = $rng.Find.Execute("Варіант 14", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "Варіант 10", 2)

# 2. Relocate the "_GoBack" bookmark: it currently sits on the "Мета роботи"
#    Heading2 paragraph; it should instead sit right after the text we just
#    changed ("Варіант 10"), as a zero-length bookmark immediately following
#    the run. Re-adding a bookmark with an existing name moves it, so the
#    old one at "Мета роботи" disappears automatically.
if ($found) {
    $target = $rng.Duplicate
    $target.Collapse(0)
    # Insert a temporary marker character so the bookmark range resolves
    # unambiguously to "right after the run, still inside this paragraph"
    # rather than sliding into the following paragraph's properties.
    $target.InsertAfter("x")
    $d.Bookmarks.Add("_GoBack", $target) | Out-Null
    $target.Text = ""
}
